$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (19-dec) before the
# existing "01-oct." column (EV), shifting 01-oct..31-oct right by one
# column (EV:FZ -> EW:GA). New column gets header "19-dec" in row 1 and
# "-" placeholders in the data rows (2-25), matching the other not-yet-
# available future dates already present in the sheet.
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("EV1:EV25").EntireColumn.Insert()
$wsPrix.Range("EV1").Value = "19-dec"
$wsPrix.Range("EV2:EV25").Value = "-"

# --- Sheet "Gaz": append the next day's price row.
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDate = $wsGaz.Range("A182")
$gazDate.Value = "'2025-12-17"
$gazDate.Style = "Normal"
$wsGaz.Range("B182").Value = 25.75

# --- Sheet "CO2": append the next day's price row.
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2Date = $wsCO2.Range("A182")
$co2Date.Value = "'2025-12-17"
$co2Date.Style = "Normal"
$wsCO2.Range("B182").Value = 84.8
